$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated cells remain stored as text (matching original inline string cells)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.31%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.19%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.431"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.73%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08112"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.96%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.720"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.28%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.07%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.42%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9460"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.25%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1183"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.85%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1892"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.85%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09701"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.02%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04204"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "8.80%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.68%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001304"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.12%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006135"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.55%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.563"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.61%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.39%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.808"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-6.55%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1362"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.01%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2610"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.16%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04390"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.49%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001244"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.23%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004319"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.29%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001242"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.88%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004024"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "32.13%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02663"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.11%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05563"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.74%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007776"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.98%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009827"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.59%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1402"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.27%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002133"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.30%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009625"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-13.01%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007119"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.05%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000757"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003489"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.76%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002290"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.58%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002121"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.04%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002020"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "1.04%"
